# Add two new columns, I (I0) and J (IF), to the existing data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell so the new headers look consistent with the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Data rows (I2:J20) --------------------------------------------------
$values = @{
    2  = @(5, 6)
    3  = @(4, 6)
    4  = @(6, 6)
    5  = @(1, 5)
    6  = @(1, 6)
    7  = @(1, 5)
    8  = @(1, 1)
    9  = @(1, 5)
    10 = @(1, 3)
    11 = @(1, 6)
    12 = @(1, 5)
    13 = @(1, 5)
    14 = @(1, 3)
    15 = @(1, 4)
    16 = @(1, 4)
    17 = @(1, 4)
    18 = @(1, 3)
    19 = @(1, 4)
    20 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

# --- Update the sheet dimension reference --------------------------------
$ws.UsedRange | Out-Null
